$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.146.53"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.893.25"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'305.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.5364"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "'0.07260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "'21.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").Value = "'0.8942"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'0.08167"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "'94.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "'5.326"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "1.810.05"
$ws.Range("E15").Value = "  -4.85%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'14.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'0.000008622"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "27.005.30"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "'5.015"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "'10.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'6.449"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'148.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'2.277"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'18.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'1.744"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'116.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'4.802"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").Value = "'0.09159"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'0.8159"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "'0.05028"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'1.210"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "'3.015"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "'3.289"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").Value = "'2.652"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").Value = "'0.5928"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "'0.01981"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "'1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'9.212"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "'6.599"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'114.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").Value = "'0.5057"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("D45").Value = "'0.1519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "'1.623"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "'37.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "'0.06069"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "'62.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.85%  "
